$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns to match the refreshed crypto data feed.
# Price cells are plain text in the source sheet (e.g. "64.583.78" uses dot as a
# thousands separator and would otherwise be re-interpreted as a number by Excel),
# so we force the Text number format before writing, then clear formats again so the
# cell keeps the worksheet default style (just like the original, un-styled cells).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.469.06"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.527.57"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.61%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.77"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.47"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.83%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.526.10"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.51%  "

$ws.Range("E10").Value = "  +0.37%  "

$ws.Range("E11").Value = "  -1.48%  "

$ws.Range("E13").Value = "  -1.10%  "

$ws.Range("E14").Value = "  +0.97%  "

$ws.Range("E15").Value = "  +0.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.990.50"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.461.93"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.527.07"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.84"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.95"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.27"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "328.31"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.34%  "

$ws.Range("E23").Value = "  +0.29%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("E25").Value = "  -0.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.36"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "647.85"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("E28").Value = "  +4.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.650.64"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.28%  "

$ws.Range("E30").Value = "  +4.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.01"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.37%  "

$ws.Range("E33").Value = "  +1.63%  "

$ws.Range("E34").Value = "  +2.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("E36").Value = "  +0.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.81"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.23%  "

$ws.Range("E38").Value = "  +2.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "153.97"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.68%  "

$ws.Range("E40").Value = "  +0.78%  "

$ws.Range("E41").Value = "  +1.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.80"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.38%  "

$ws.Range("E43").Value = "  +2.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "163.15"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +6.85%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0300"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.52"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.74%  "

$ws.Range("E48").Value = "  +1.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.24"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.51%  "

$ws.Range("E50").Value = "  +1.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0518"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.32%  "
